# Updates the crypto price table (rows 2-51) with refreshed Price (D) and
# Volume(1h) (E) figures, including a couple of Coin/Link swaps (rows 26-27,
# ImmutableX <-> LidoDAOToken) from the latest GitHub Actions refresh.
#
# Values that look like plain decimals (e.g. "21.80", "1.002") are written
# through a text-formatted cell so Excel keeps them as literal strings
# (preserving trailing zeros / leading dots) instead of silently coercing
# them to numbers; everything else is assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.238.95"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").Value = "1.786.00"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.02"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3836"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -3.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.45"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -3.42%  "
$ws.Range("E10").Value = "  -3.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07482"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -3.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.80"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E14").Value = "  -2.82%  "
$ws.Range("D15").Value = "1.785.85"
$ws.Range("E15").Value = "  -1.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.069"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06653"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.72"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -3.66%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.604"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.32"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D23").Value = "27.241.30"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.35"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -6.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.400"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -2.74%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.534"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -5.89%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.479"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -1.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.23"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -3.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.50"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("D30").Value = "1.988.26"
$ws.Range("E30").Value = "  -1.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "134.16"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.008"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.063"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -5.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08721"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -1.31%  "
$ws.Range("E35").Value = "  -4.90%  "
$ws.Range("E36").Value = "  -4.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6904"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -2.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.435"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -4.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2200"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -3.14%  "
$ws.Range("E40").Value = "  -3.53%  "
$ws.Range("E41").Value = "  -3.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.778"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -2.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.238"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -4.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.27"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -5.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6487"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -1.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.852"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -3.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.139"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -2.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.86"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -2.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07139"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -3.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.78"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -2.68%  "
